$d = $word.ActiveDocument

# --- Step 1: Clear paragraph 2's run ("When should you use Redux?" -> empty paragraph) ---
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$r2.MoveEnd(1, -1)
$r2.Delete()

# --- Step 2: Replace paragraph 3's text with the full TypeScript answer (multi-run) ---
$p3 = $d.Paragraphs(3)
$p3xml = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0"/><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t>TypeScript is a</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t>n open-source</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t xml:space="preserve"> JS based coding </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t xml:space="preserve">created by Microsoft. It is </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t>a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t xml:space="preserve"> object oriented, compiled language.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/></w:rPr><w:t xml:space="preserve"> It e</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t>xtends JavaScript and improves the developer experience. It enables developers to add type safety to their projects. TypeScript</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t xml:space="preserve"> also provides more functionality to</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t xml:space="preserve"> various other features, like interfaces, type aliases, abstract classes, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t>function</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t xml:space="preserve"> overloading, tuple, generics, etc.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t xml:space="preserve"> It can also help build and manage large-scale JS objects which are also reasons many companies use it.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p3.Range.InsertXML($p3xml)

# --- Step 3: Append the new paragraphs after the (now TypeScript-answer) paragraph 3.
#     Done in two batches of 3 paragraphs each (InsertXML has issues with 4+ paragraphs
#     in a single call), targeting the still-empty trailing paragraph's range so the
#     document keeps a single blank paragraph at the very end, matching the target. ---
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$batch1 = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0"/><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t>https://www.contentful.com/</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0"/><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0"/><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t>When should you use Redux?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$pLast.Range.InsertXML($batch1)

$pLast2 = $d.Paragraphs($d.Paragraphs.Count)
$batch2 = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0"/><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t>When you have large amounts of state in an application that is needed in multiple points through out the app. If the state needs to be updated frequently. When the logic to a state may be complex. If the app has a medium-large sized codebase and/or many individuals are working on the app. If you need to see how the state is being updated over time.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0"/><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0"/><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="21252A"/></w:rPr><w:t>https://redux.js.org/</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$pLast2.Range.InsertXML($batch2)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
foreach ($p in $d.Paragraphs) {
    Write-Output ("PARA: [" + $p.Range.Text + "]")
}
